{"js": "const replacements = [\n  [\"817\u00f75=\", \"561\u00f78=\"],\n  [\"772\u00f75=\", \"676\u00f74=\"],\n  [\"353\u00f76=\", \"751\u00f76=\"],\n  [\"303\u00f79=\", \"933\u00f79=\"],\n  [\"530\u00f73=\", \"152\u00f72=\"],\n  [\"854\u00f73=\", \"753\u00f72=\"],\n  [\"576\u00f78=\", \"710\u00f76=\"],\n  [\"584\u00f72=\", \"958\u00f77=\"],\n  [\"975\u00f74=\", \"781\u00f73=\"],\n  [\"312\u00f78=\", \"540\u00f72=\"],\n  [\"368\u00f77=\", \"495\u00f77=\"],\n  [\"987\u00f73=\", \"836\u00f79=\"],\n  [\"645\u00f78=\", \"557\u00f78=\"],\n  [\"235\u00f78=\", \"285\u00f74=\"],\n  [\"679\u00f74=\", \"266\u00f79=\"],\n  [\"914\u00f79=\", \"633\u00f75=\"],\n  [\"705\u00f75=\", \"182\u00f79=\"],\n  [\"196\u00f72=\", \"584\u00f79=\"],\n  [\"297\u00f77=\", \"461\u00f78=\"],\n  [\"602\u00f74=\", \"485\u00f76=\"],\n  [\"921\u00f79=\", \"518\u00f74=\"],\n  [\"348\u00f79=\", \"553\u00f77=\"],\n  [\"462\u00f72=\", \"627\u00f73=\"],\n  [\"204\u00f75=\", \"371\u00f72=\"],\n  [\"245\u00f72=\", \"472\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each problem's \"NNN\u00f7N=\" label is unique in the document, so a simple\n# whole-content Find/Replace (MatchCase, no wildcards) per pair is safe.\n$replacements = @(\n    @(\"817\u00f75=\", \"561\u00f78=\"),\n    @(\"772\u00f75=\", \"676\u00f74=\"),\n    @(\"353\u00f76=\", \"751\u00f76=\"),\n    @(\"303\u00f79=\", \"933\u00f79=\"),\n    @(\"530\u00f73=\", \"152\u00f72=\"),\n    @(\"854\u00f73=\", \"753\u00f72=\"),\n    @(\"576\u00f78=\", \"710\u00f76=\"),\n    @(\"584\u00f72=\", \"958\u00f77=\"),\n    @(\"975\u00f74=\", \"781\u00f73=\"),\n    @(\"312\u00f78=\", \"540\u00f72=\"),\n    @(\"368\u00f77=\", \"495\u00f77=\"),\n    @(\"987\u00f73=\", \"836\u00f79=\"),\n    @(\"645\u00f78=\", \"557\u00f78=\"),\n    @(\"235\u00f78=\", \"285\u00f74=\"),\n    @(\"679\u00f74=\", \"266\u00f79=\"),\n    @(\"914\u00f79=\", \"633\u00f75=\"),\n    @(\"705\u00f75=\", \"182\u00f79=\"),\n    @(\"196\u00f72=\", \"584\u00f79=\"),\n    @(\"297\u00f77=\", \"461\u00f78=\"),\n    @(\"602\u00f74=\", \"485\u00f76=\"),\n    @(\"921\u00f79=\", \"518\u00f74=\"),\n    @(\"348\u00f79=\", \"553\u00f77=\"),\n    @(\"462\u00f72=\", \"627\u00f73=\"),\n    @(\"204\u00f75=\", \"371\u00f72=\"),\n    @(\"245\u00f72=\", \"472\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format,\n    # ReplaceWith, Replace(2=wdReplaceAll)\n    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
